# Add a new "BPFUbIP" acronym row (indst / BAU Percentage Fuel Use by
# Industrial Process) to the "Key to Variables" sheet, just above the
# existing "indst / BPoIFUfE" row (i.e. it becomes the new row 162).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")
$ws.Activate()

# Insert a brand-new row at position 162; this pushes the former rows
# 162-275 down to 163-276 and copies row 161's formatting by default.
$ws.Rows.Item(162).Insert()

# Populate the new row's contents.
$ws.Cells.Item(162, 1).Value = "indst"
$ws.Cells.Item(162, 2).Value = "BPFUbIP"
$ws.Cells.Item(162, 3).Value = "BAU Percentage Fuel Use by Industrial Process"
$ws.Cells.Item(162, 6).Value = "low"

# The inserted row picked up row 161's formatting (which includes a
# populated column D and a different "Importance" fill colour). Clear the
# stray column D cell and recolor column F to match the "low" style used
# elsewhere in the sheet (copy formats only, from a reference "low" cell).
$ws.Cells.Item(162, 4).Clear()

$ws.Cells.Item(5, 6).Copy()
$ws.Cells.Item(162, 6).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the view: header row frozen, scrolled down near the new row, with
# the new row's acronym cell selected (matches the saved sheet view state).
[void]($excel.ActiveWindow.FreezePanes = $false)
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
$excel.ActiveWindow.ScrollRow = 155
[void]$ws.Range("A162").Select()
